$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Next empty row after existing data (row 52 has data -> new row is 53)
$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "G6"
$ws.Cells.Item($newRow, 2).Value = "Spend 10 Hours without phone"

# Date column: copy the style/number format used by the Date column above (row 52)
$ws.Cells.Item($newRow - 1, 3).Copy() | Out-Null
$ws.Cells.Item($newRow, 3).PasteSpecial(-4122) | Out-Null ## xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item($newRow, 3).Value = 45907

$ws.Cells.Item($newRow, 4).Value = 1
$ws.Cells.Item($newRow, 5).Value = 0
$ws.Cells.Item($newRow, 6).Value = 0
